$wb = $excel.ActiveWorkbook
$wsProdutos  = $wb.Worksheets.Item("produtos")
$wsMovimentos = $wb.Worksheets.Item("movimentos")

# --- "produtos" sheet: three new product rows (90-92) ---------------------
$wsProdutos.Range("A90").Value() = 89
$wsProdutos.Range("B90").Value() = "vassoura 3"
# Quantity column holds a digit-looking label, not a real number - force
# text storage the same way a user would in Excel (leading apostrophe).
$wsProdutos.Range("D90").Value() = "'4"
$wsProdutos.Range("E90").Value() = 0

$wsProdutos.Range("A91").Value() = 90
$wsProdutos.Range("B91").Value() = "'44"
$wsProdutos.Range("D91").Value() = "kg"
$wsProdutos.Range("E91").Value() = 0

$wsProdutos.Range("A92").Value() = 91
$wsProdutos.Range("B92").Value() = "arroz 5L"
$wsProdutos.Range("C92").Value() = ""
$wsProdutos.Range("D92").Value() = "5 L"
$wsProdutos.Range("E92").Value() = 0

# --- "movimentos" sheet: four new finished movement rows (11-14) ----------
$wsMovimentos.Range("A11").Value() = 10
$wsMovimentos.Range("B11").Value() = 1
$wsMovimentos.Range("C11").Value() = "ENTRADA"
$wsMovimentos.Range("D11").Value() = 11
$wsMovimentos.Range("E11").Value() = "2025-12-05 16:27:50"

$wsMovimentos.Range("A12").Value() = 11
$wsMovimentos.Range("B12").Value() = 49
$wsMovimentos.Range("C12").Value() = "ENTRADA"
$wsMovimentos.Range("D12").Value() = 11
$wsMovimentos.Range("E12").Value() = "2025-12-09 14:44:36"

$wsMovimentos.Range("A13").Value() = 12
$wsMovimentos.Range("B13").Value() = 1
$wsMovimentos.Range("C13").Value() = "ENTRADA"
$wsMovimentos.Range("D13").Value() = 11
$wsMovimentos.Range("E13").Value() = "2025-12-09 15:22:50"

$wsMovimentos.Range("A14").Value() = 13
$wsMovimentos.Range("B14").Value() = 84
$wsMovimentos.Range("C14").Value() = "ENTRADA"
$wsMovimentos.Range("D14").Value() = 22
$wsMovimentos.Range("E14").Value() = "2025-12-09 15:22:50"
